$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("instructions")

$newText = "Sprawia mi radość, jeśli mam jasne preferencje w przypadku nowych rzeczy.`n`n1 = bardzo nietypowe/niecharakterystyczne`n2 = nieco nietypowe/niecharakterystyczne`n3 = ani charakterystyczny, ani nietypowy`n4 = nieco charakterystyczne`n5 = bardzo charakterystyczne`n`nKliknij linię, a następnie potwierdź swój wybór, klikając szary przycisk poniżej."

$ws.Range("A8").Value = $newText
